# Add automatic publist, add matomo
#
# 1. Rename the existing "Sheet1" to "Published".
# 2. Add a new "Preprints" worksheet right after it.
# 3. Populate the new sheet's header row (copied format+values from Published!A1:H1)
#    and a first data row for the new preprint, including live hyperlinks.
# 4. Re-point the UI selection: Published shows a selected header row (no longer the
#    active tab), Preprints becomes the active tab with G2 selected.

$wb = $excel.ActiveWorkbook

$published = $wb.Worksheets.Item(1)
$published.Name = "Published"

# New sheet goes directly after "Published".
$preprints = $wb.Worksheets.Add($null, $published)
$preprints.Name = "Preprints"

# --- column widths (mirrors Published, minus column B which stays default) ---
$preprints.Columns.Item(1).ColumnWidth = 38.83203125
$preprints.Columns.Item(3).ColumnWidth = 12.33203125
$preprints.Columns.Item(4).ColumnWidth = 91.83203125
$preprints.Columns.Item(5).ColumnWidth = 41.83203125
$preprints.Columns.Item(6).ColumnWidth = 57.83203125
$preprints.Columns.Item(7).ColumnWidth = 53.83203125
$preprints.Columns.Item(8).ColumnWidth = 36.6640625

# --- header row: same labels/style as Published!A1:H1 ---
$published.Range("A1:H1").Copy($preprints.Range("A1:H1"))

# --- first data row ---
$preprints.Range("B2").Value = 2025
$preprints.Range("C2").Value = "true"
$preprints.Range("D2").Value = "Interpersonal versus intrapersonal emotion regulation: Intensity of negative emotion predicts usage probability"
$preprints.Range("G2").Value = "https://osf.io/hjzpw/"
$preprints.Range("H2").Value = "https://osf.io/dwnya"

$preprints.Range("A2").Value = "https://doi.org/10.31234/osf.io/4u8kj_v1"
$preprints.Range("E2").Value = "https://doi.org/10.31234/osf.io/4u8kj_v1"

$preprints.Hyperlinks.Add($preprints.Range("A2"), "https://doi.org/10.31234/osf.io/4u8kj_v1") | Out-Null
$preprints.Hyperlinks.Add($preprints.Range("E2"), "https://doi.org/10.31234/osf.io/4u8kj_v1") | Out-Null

# --- page margins: Preprints uses the (Mac-default) 2 cm top/bottom margins ---
$preprints.PageSetup.TopMargin = 56.692913399999995
$preprints.PageSetup.BottomMargin = 56.692913399999995

# --- selection / active-tab bookkeeping ---
$published.Rows.Item(1).Select()

$preprints.Activate()
$preprints.Range("G2").Select()

Write-Output "ok"
